# Commit: "add DNA data and relabel nadp sheet"
$wb = $excel.ActiveWorkbook

$wsNadp = $wb.Worksheets.Item(1)
$wsDna  = $wb.Worksheets.Item(2)

# --- Relabel the "nad" sheet to "nadp" ---
$wsNadp.Name = "nadp"

# --- New DNA quantity columns (F:H) for rows 2-25 on the "dna" sheet ---
# Rows 2-9 (the standard-curve block) get an explicit black font; rows
# 10-25 (the sample rows) keep the default style, same as the source data.
$rows = @(
  @(2, 45479, 47505, 47532),
  @(3, 594246, 592748, 607592),
  @(4, 1094433, 1084946, 1085169),
  @(5, 2116238, 2083453, 2167345),
  @(6, 3969176, 4062854, 4039748),
  @(7, 7683769, 7525415, 7405485),
  @(8, 14950687, 15095009, 15088696),
  @(9, 32258614, 31554882, 32605324),
  @(10, 7110940, 7055460, 7593087),
  @(11, 6999632, 7918300, 7455102),
  @(12, 6945892, 7045710, 7908128),
  @(13, 6540020, 6804756, 6497860),
  @(14, 6447102, 6393064, 6894206),
  @(15, 5488896, 5820153, 5711590),
  @(16, 2610636, 3904350, 3136986),
  @(17, 2710910, 3625672, 3231936),
  @(18, 6466554, 5465786, 4856718),
  @(19, 5479005, 4911644, 4144374),
  @(20, 4634501, 5218002, 4121988),
  @(21, 3217068, 3201932, 4213901),
  @(22, 3417824, 3797327, 3236893),
  @(23, 2583024, 3146347, 3609966),
  @(24, 2580380, 3082436, 2653759),
  @(25, 1733593, 1932103, 2108261)
)

foreach ($row in $rows) {
  $r = $row[0]
  $wsDna.Cells.Item($r, 6).Value = $row[1]
  $wsDna.Cells.Item($r, 7).Value = $row[2]
  $wsDna.Cells.Item($r, 8).Value = $row[3]
}

$wsDna.Range("F2:H9").Font.Color = 0

# --- Selections: "dna" sheet remembers F18:H25, then "nadp" is reactivated
#     (and keeps the selected/active tab) with C30 selected. ---
$wsDna.Activate()
$wsDna.Range("F18:H25").Select()

$wsNadp.Activate()
$wsNadp.Range("C30").Select()
